# Fruta / hortaliza, semanal
# Insert two new weekly price records at rows 309-310 (pushing the
# existing rows 309..376 down to 311..378).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 309, shifting the rest down.
$ws.Rows("309:310").Insert()

# --- New row 309: Papa / Asterix / 1a (guarda) ---
$ws.Cells.Item(309, 1).Value = 7
$ws.Cells.Item(309, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(309, 3).Value = "Ñuble"
$ws.Cells.Item(309, 4).Value = 44785
$ws.Cells.Item(309, 5).Value = 16
$ws.Cells.Item(309, 6).Value = 100114001
$ws.Cells.Item(309, 7).Value = "Papa"
$ws.Cells.Item(309, 8).Value = "Asterix"
$ws.Cells.Item(309, 9).Value = "1a (guarda)"
$ws.Cells.Item(309, 10).Value = 120
$ws.Cells.Item(309, 11).Value = 7000
$ws.Cells.Item(309, 12).Value = 7500
$ws.Cells.Item(309, 13).Value = 7250
$ws.Cells.Item(309, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(309, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(309, 16).Value = 290
$ws.Cells.Item(309, 17).Value = 25
$ws.Cells.Item(309, 18).Value = "Hortaliza"

# --- New row 310: Papa / Patagonia / 1a (guarda) ---
$ws.Cells.Item(310, 1).Value = 7
$ws.Cells.Item(310, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(310, 3).Value = "Ñuble"
$ws.Cells.Item(310, 4).Value = 44785
$ws.Cells.Item(310, 5).Value = 16
$ws.Cells.Item(310, 6).Value = 100114001
$ws.Cells.Item(310, 7).Value = "Papa"
$ws.Cells.Item(310, 8).Value = "Patagonia"
$ws.Cells.Item(310, 9).Value = "1a (guarda)"
$ws.Cells.Item(310, 10).Value = 120
$ws.Cells.Item(310, 11).Value = 7000
$ws.Cells.Item(310, 12).Value = 7500
$ws.Cells.Item(310, 13).Value = 7250
$ws.Cells.Item(310, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(310, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(310, 16).Value = 290
$ws.Cells.Item(310, 17).Value = 25
$ws.Cells.Item(310, 18).Value = "Hortaliza"
